$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.330.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "'3.777.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'625.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").Value = "'166.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("D7").Value = "'3.776.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("D12").Value = "'6.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "'35.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "'4.416.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'3.778.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'69.331.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "'17.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.114"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'7.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'469.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").Value = "'9.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "'0.706"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").Value = "'0.0000148"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.46%  "
$ws.Range("D25").Value = "'83.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'12.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("E27").Value = "  +4.38%  "
$ws.Range("D28").Value = "'10.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'3.927.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +3.89%  "
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").Value = "'7.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'28.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "'3.729.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'9.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  +11.90%  "
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").Value = "'3.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.52%  "
$ws.Range("D41").Value = "'5.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").Value = "'0.967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D45").Value = "'0.298"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'153.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'43.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "'46.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("D50").Value = "'8.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").Value = "'1.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
